$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 90.33701601336995;   C = 97.99789108864614;    D = 99.25719803621033;    E = 98.72966839909351;  F = 98.23074575766434;  G = 97.28828811886578;  H = 95.86270046195699 }
    3  = @{ B = 76.89876022041162;   C = 94.07007365722883;    D = 99.00016420823033;    E = 98.71570471705476;  F = 98.45910039451762;  G = 97.55691708871703;  H = 96.12898886077714 }
    4  = @{ B = 91.92264418048704;   C = 98.13284670726674;    D = 99.23778052608313;    E = 98.71004261302667;  F = 98.21268306562727;  G = 97.2763229845675;   H = 95.85348005711388 }
    5  = @{ B = 78.53890916540075;   C = 97.06733670160607;    D = 99.41478311448603;    E = 98.8559289049183;   F = 98.30016816553849;  G = 97.32977735399146;  H = 95.88960392010006 }
    6  = @{ B = 85.39860290829367;   C = 95.47603018624659;    D = 99.03557576777675;    E = 98.57199874879504;  F = 98.31066284644379;  G = 97.56666974420268;  H = 96.2247806100732  }
    7  = @{ B = 41.39579404654932;   C = 50.28889861122961;    D = 98.82604622396775;    E = 98.92103190527189;  F = 98.49091268988114;  G = 97.59349024792978;  H = 96.19316105527326 }
    8  = @{ B = 79.57218880010741;   C = 86.35384427366763;    D = 99.32680854529019;    E = 98.92216067459138;  F = 98.47245700686243;  G = 97.56316972844064;  H = 96.13267281987892 }
    9  = @{ B = 74.61769130626827;   C = 95.16083180530659;    D = 99.37025810399646;    E = 98.88526414555285;  F = 98.40341225926801;  G = 97.4590900614585;   H = 96.02337987804562 }
    10 = @{ B = 0.1688198682658943;  C = 0.04790827893590049;  D = -0.001756253530937295; E = 0.01675237852423847; F = 18.26692517429614; G = 98.11386914362433; H = 95.52771486509967 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Range("$col$row").Value = $values[$row][$col]
    }
}
